# ParaBank-AutomationSuite/testdata/registerdata.xlsx
# Rename the placeholder "Sheet1" tab to "billPageDataTest" and replace its
# sample bill-pay data with a new set of test values (plus matching column
# widths / selection), per the commit "add new test xml, test class and
# update the excel data driven file".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "billPageDataTest"

# ---- Row 1: header labels (reordered so "Name" leads, "Phone"/"Account"/
# "verify Account" shift before "Amount"/"From Acount") ----
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Address"
$ws.Range("C1").Value = "City"
$ws.Range("D1").Value = "State"
$ws.Range("E1").Value = "ZipCode"
$ws.Range("F1").Value = "Phone"
$ws.Range("G1").Value = "Account"
$ws.Range("H1").Value = "verify Account"
$ws.Range("I1").Value = "Amount"
$ws.Range("J1").Value = "From Acount"

# ---- Row 2: sample data ----
# Format A2:H2 as Text *before* writing the values so numeric-looking
# strings (zip code, phone number, account numbers with a leading zero)
# are stored as genuine text rather than being coerced to numbers.
$ws.Range("A2:H2").NumberFormat = "@"
$ws.Range("B2").Value = "Iyara Ijumu LGA Kogi State"
$ws.Range("C2").Value = "Kabba"
$ws.Range("D2").Value = "Kogi"
$ws.Range("E2").Value = "49292"
$ws.Range("F2").Value = "09188282828"
$ws.Range("G2").Value = "20004"
$ws.Range("H2").Value = "20004"
$ws.Range("A2").Value = "Alemidan Ojo"

# I2 keeps a genuine numeric value but still picks up the Text cell style;
# set the value first, then apply the format so the stored <v> stays numeric.
$ws.Range("I2").Value = 100
$ws.Range("I2").NumberFormat = "@"

# J2 is a plain, unstyled number.
$ws.Range("J2").Value = 1

# ---- Column widths (best fit to content) ----
$ws.Columns.Item(1).ColumnWidth = 12.140625
$ws.Columns.Item(2).ColumnWidth = 8.140625
$ws.Columns.Item(3).ColumnWidth = 5.85546875
$ws.Columns.Item(4).ColumnWidth = 5.85546875
$ws.Columns.Item(5).ColumnWidth = 8.28515625
$ws.Columns.Item(7).ColumnWidth = 8.140625
$ws.Columns.Item(8).ColumnWidth = 13.85546875
$ws.Columns.Item(9).ColumnWidth = 8.140625
$ws.Columns.Item(10).ColumnWidth = 12.28515625

# ---- Selection / scroll position ----
$ws.Activate() | Out-Null
$ws.Range("J2").Select() | Out-Null
